# Update cryptos list values as scraped on Sat Apr  6 06:57:45 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.056.23"
$ws.Cells.Item(2, 5).Value = "  +1.59%  "
$ws.Cells.Item(3, 4).Value = "3.343.09"
$ws.Cells.Item(3, 5).Value = "  +1.95%  "
$cell = $ws.Cells.Item(4, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "582.26"
$cell.Style = $origStyle
$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "177.61"
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  +1.51%  "
$cell = $ws.Cells.Item(7, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = $origStyle
$ws.Cells.Item(7, 5).Value = "  -0.27%  "
$ws.Cells.Item(8, 5).Value = "  +1.60%  "
$ws.Cells.Item(9, 4).Value = "3.338.74"
$ws.Cells.Item(9, 5).Value = "  +2.03%  "
$ws.Cells.Item(10, 5).Value = "  +6.07%  "
$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.583"
$cell.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  +2.06%  "
$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "47.29"
$cell.Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  +3.98%  "
$ws.Cells.Item(13, 5).Value = "  +2.32%  "
$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "686.50"
$cell.Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  -0.34%  "
$ws.Cells.Item(15, 4).Value = "3.882.43"
$ws.Cells.Item(15, 5).Value = "  +2.05%  "
$ws.Cells.Item(16, 5).Value = "  +2.24%  "
$ws.Cells.Item(17, 4).Value = "68.105.43"
$ws.Cells.Item(17, 5).Value = "  +1.62%  "
$ws.Cells.Item(18, 5).Value = "  -0.23%  "
$ws.Cells.Item(19, 4).Value = "3.341.93"
$ws.Cells.Item(19, 5).Value = "  +1.80%  "
$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.47"
$cell.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  +1.03%  "
$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.12"
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  +3.84%  "
$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.900"
$cell.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  +1.67%  "
$ws.Cells.Item(23, 5).Value = "  +5.03%  "
$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.09"
$cell.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  +1.16%  "
$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "99.27"
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  +0.16%  "
$ws.Cells.Item(26, 5).Value = "  +1.01%  "
$ws.Cells.Item(27, 5).Value = "  +0.59%  "
$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.58"
$cell.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  +3.57%  "
$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "33.05"
$cell.Style = $origStyle
$ws.Cells.Item(29, 5).Value = "  +0.29%  "
$cell = $ws.Cells.Item(30, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.60"
$cell.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  +2.75%  "
$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.15"
$cell.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  +6.22%  "
$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "571.63"
$cell.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  -0.15%  "
$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.03"
$cell.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  +2.22%  "
$ws.Cells.Item(34, 5).Value = "  +2.94%  "
$ws.Cells.Item(35, 4).Value = "3.723.51"
$ws.Cells.Item(35, 5).Value = "  -4.21%  "
$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "57.31"
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  +3.26%  "
$ws.Cells.Item(37, 5).Value = "  +0.09%  "
$ws.Cells.Item(38, 5).Value = "  +0.86%  "
$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "34.75"
$cell.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  +9.65%  "
$ws.Cells.Item(40, 5).Value = "  +2.99%  "
$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.22"
$cell.Style = $origStyle
$ws.Cells.Item(41, 5).Value = "  +7.22%  "
$ws.Cells.Item(42, 5).Value = "  +2.83%  "
$ws.Cells.Item(43, 4).Value = "0.0₃0680"
$ws.Cells.Item(43, 5).Value = "  +1.60%  "
$ws.Cells.Item(44, 2).Value = "ApeXProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.35"
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  -0.26%  "
$ws.Cells.Item(45, 2).Value = "TheGraph"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.338"
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +3.45%  "
$ws.Cells.Item(46, 5).Value = "  +0.90%  "
$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.Style = $origStyle
$ws.Cells.Item(47, 5).Value = "  +6.12%  "
$ws.Cells.Item(48, 5).Value = "  +1.34%  "
$ws.Cells.Item(49, 5).Value = "  -0.35%  "
$ws.Cells.Item(50, 5).Value = "  -2.80%  "
$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "129.94"
$cell.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  -0.53%  "
